# SummonTable.xlsx edit: rework Skills sheet summon data (W/A items -> S0001-S0006),
# shrink table from 62 to 32 rows, update JSON export paths, switch active tab to Skills.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Follower sheet: keep its JSON path string's *value* identical, but touch
#    it first so the shared-string table keeps SummonTableFollower.json where
#    SummonTableSkills.json used to sit (matches upstream reorder), and the
#    new Skills path + new S000x strings get appended after it in the order
#    they are created below.
# ---------------------------------------------------------------------------
$wsFollower = $wb.Worksheets.Item("Follower")
$wsFollower.Range("A1").Value = "D:\Project\TeamProject-IdleGame\IdleGame\Assets\Resources\Texts\SummonTableFollower.json"

# ---------------------------------------------------------------------------
# 2) Skills sheet: rebuild the summon-item rows.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Skills")

# -- Grade 1 "base" row (3-8): same probabilities, new item ids --
$ws.Range("B3").Value = "S0001"
$ws.Range("B4").Value = "S0002"
$ws.Range("B5").Value = "S0003"
$ws.Range("B6").Value = "S0004"
$ws.Range("B7").Value = "S0005"
$ws.Range("B8").Value = "S0006"

# -- Grade 2 (rows 9-14): grade number, then items mirroring rows 3-8 --
$ws.Range("A9").Value = 2
$ws.Range("A10").Formula = "=A9"
$ws.Range("A11:A14").Formula = "=A10"
$ws.Range("B9").Formula = "=B3"
$ws.Range("B10").Formula = "=B4"
$ws.Range("B11").Formula = "=B5"
$ws.Range("B12").Formula = "=B6"
$ws.Range("B13").Formula = "=B7"
$ws.Range("B14").Formula = "=B8"
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 4995
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("C14").Value = 1

# -- Grade 3 (rows 15-20) --
$ws.Range("A15").Value = 3
$ws.Range("A16").Formula = "=A15"
$ws.Range("A17:A20").Formula = "=A16"
$ws.Range("B15").Formula = "=B9"
$ws.Range("B16").Formula = "=B10"
$ws.Range("B17").Formula = "=B11"
$ws.Range("B18").Formula = "=B12"
$ws.Range("B19").Formula = "=B13"
$ws.Range("B20").Formula = "=B14"
$ws.Range("C15").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("C17").Value = 4995
$ws.Range("C18").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("C20").Value = 1

# -- Grade 4 (rows 21-26) --
$ws.Range("A21").Value = 4
$ws.Range("A22").Formula = "=A21"
$ws.Range("A23:A26").Formula = "=A22"
$ws.Range("B21").Formula = "=B15"
$ws.Range("B22").Formula = "=B16"
$ws.Range("B23").Formula = "=B17"
$ws.Range("B24").Formula = "=B18"
$ws.Range("B25").Formula = "=B19"
$ws.Range("B26").Formula = "=B20"
$ws.Range("C21").Value = 1
$ws.Range("C22").Value = 1
$ws.Range("C23").Value = 1
$ws.Range("C24").Value = 4995
$ws.Range("C25").Value = 1
$ws.Range("C26").Value = 1

# -- Grade 5 (rows 27-32) --
$ws.Range("A27").Value = 5
$ws.Range("A28").Formula = "=A27"
$ws.Range("A29:A32").Formula = "=A28"
$ws.Range("B27").Formula = "=B21"
$ws.Range("B28").Formula = "=B22"
$ws.Range("B29").Formula = "=B23"
$ws.Range("B30").Formula = "=B24"
$ws.Range("B31").Formula = "=B25"
$ws.Range("B32").Formula = "=B26"
$ws.Range("C27").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("C30").Value = 1
$ws.Range("C31").Value = 4995
$ws.Range("C32").Value = 1

# -- Drop the now-unused tail (old rows 33-62 plus the D-column SUM helpers) --
$ws.Range("A33:D62").Clear()
$ws.Range("D1:D32").Clear()

# -- Resize the Skills table (표1_3) down to the new extent --
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A2:C32"))

# -- Point the export-path cell at the new Skills JSON location --
$ws.Range("A1").Value = "D:\Project\TeamProject-IdleGame\IdleGame\Assets\Resources\Texts\Summon\SummonTableSkills.json"

# ---------------------------------------------------------------------------
# 3) View state: Skills tab becomes the active/selected tab instead of Follower.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E8").Select()
